$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters being updated per row: B,C,D,E,F,I,J,K,L,M,N (column H is intentionally blank)
$cols = @(2,3,4,5,6,9,10,11,12,13,14)

# Each row: @(rowNumber, <11 new values for cols B,C,D,E,F,I,J,K,L,M,N>)
$data = @(
    @(2, 1.02, 1.026776104332363, 1.028168502927843, 1.035590835825034, 1.043334832208573, 1.030869092779453, 1.031937662606526, 1.030986104273847, 1.038387003017763, 1.046108989219498, 1.014683953469651),
    @(3, 1.02, 1.027662540191614, 1.028925906230195, 1.036412375578092, 1.044319773191476, 1.030973917595845, 1.032464326535692, 1.031551522672326, 1.03901793702601, 1.046904489354026, 1.014857866084385),
    @(4, 1.02, 1.028236112692035, 1.029416288544316, 1.036944423163843, 1.044957871108046, 1.03103951742363, 1.032804479084241, 1.031917017590167, 1.039426001856217, 1.047419391362872, 1.014970176291242),
    @(5, 1.02, 1.02847723873791, 1.029622513552302, 1.037168203734273, 1.045226311604224, 1.03106656107481, 1.032947326204013, 1.032070582174762, 1.039597505188587, 1.047635893178529, 1.015017337733346),
    @(6, 1.02, 1.028517724600674, 1.029657143597271, 1.037205783752251, 1.045271394699723, 1.031071070451358, 1.032971301856974, 1.032096361061044, 1.039626298553353, 1.047672246917082, 1.015025253181992),
    @(7, 1.02, 1.028239334646708, 1.029419043866635, 1.036947412909239, 1.044961457301429, 1.031039880883689, 1.032806388417171, 1.031919069880823, 1.039428293678865, 1.047422284127007, 1.014970806676534),
    @(8, 1.02, 1.027075680960212, 1.02842440958092, 1.035868383621116, 1.043667536229189, 1.030904979757577, 1.032115781764414, 1.03117726564759, 1.0386002690625, 1.04637779822936, 1.014742773955992),
    @(9, 1.02, 1.025025146851316, 1.026674039558495, 1.033970557920569, 1.041393483544877, 1.030650236956271, 1.030894038607636, 1.029867342990371, 1.03713977024284, 1.044538566887962, 1.014339262093099),
    @(10, 1.02, 1.023658177755056, 1.025508765462042, 1.03270782339755, 1.039881568555668, 1.030469012166293, 1.030076380106953, 1.028992271578312, 1.036165225916339, 1.043313355827564, 1.014069142382432),
    @(11, 1.02, 1.023066292441461, 1.025004596599121, 1.032161652424824, 1.039227887387438, 1.030387847797206, 1.029721588741837, 1.028612946359686, 1.03574304394319, 1.042783065730629, 1.013951918773756),
    @(12, 1.02, 1.022846444205076, 1.024817387798996, 1.031958871993114, 1.038985230896155, 1.030357296004732, 1.029589693366515, 1.028471987115173, 1.03558619819226, 1.04258612884434, 1.013908338027328),
    @(13, 1.02, 1.022893602184549, 1.024857541894281, 1.032002364907666, 1.039037274736651, 1.030363867724315, 1.029617990335558, 1.02850222609207, 1.035619843412548, 1.042628370804687, 1.013917687989565),
    @(14, 1.02, 1.023048119623023, 1.02498912058867, 1.032144888650428, 1.039207826262432, 1.030385330607714, 1.029710688468005, 1.028601295861636, 1.035730079600615, 1.042766786115027, 1.013948317162497),
    @(15, 1.02, 1.023143323522779, 1.025070198820172, 1.032232714437417, 1.039312928577012, 1.030398501125004, 1.02976778829451, 1.028662327971934, 1.035797996003211, 1.042852073225647, 1.013967183692897),
    @(16, 1.02, 1.023697459778919, 1.025542234096935, 1.032744083759534, 1.039924972154514, 1.030474342114769, 1.030099910946799, 1.029017437530386, 1.036193240673858, 1.043348554480093, 1.014076916682218),
    @(17, 1.02, 1.024045061313368, 1.025838438173849, 1.033065014251713, 1.040309156400789, 1.030521194808152, 1.030308045406118, 1.029240078482861, 1.036441115062044, 1.043660047691028, 1.014145679950725),
    @(18, 1.02, 1.024247813390369, 1.026011247737463, 1.033252265509735, 1.04053333980819, 1.030548263316022, 1.030429375258715, 1.029369901207479, 1.036585676832283, 1.043841758941164, 1.014185763290992),
    @(19, 1.02, 1.024316946911256, 1.026070177857817, 1.033316123132282, 1.040609796668866, 1.030557448872702, 1.030470733449611, 1.029414160585647, 1.036634965373848, 1.043903721564889, 1.014199426403277),
    @(20, 1.02, 1.024007766745137, 1.025806654262087, 1.033030575444633, 1.040267927206386, 1.030516194831614, 1.030285721931542, 1.029216195342225, 1.036414522460154, 1.043626625077915, 1.014138304897344),
    @(21, 1.02, 1.023002617960348, 1.024950372220252, 1.032102916415299, 1.039157598942838, 1.030379021461806, 1.029683394208085, 1.028572123957522, 1.03569761857017, 1.042726025239333, 1.013939298699149),
    @(22, 1.02, 1.02237066658569, 1.024412352949612, 1.031520192112748, 1.038460358861571, 1.030290439658676, 1.029304050591385, 1.028166818535863, 1.03524670792017, 1.042159993979074, 1.013813952007541),
    @(23, 1.02, 1.022705673171196, 1.02469753259767, 1.031829054411483, 1.038829896350374, 1.030337619646361, 1.029505207770061, 1.028381711606214, 1.035485759416504, 1.042460037384964, 1.013880421727869),
    @(24, 1.02, 1.024024618549133, 1.025821015924145, 1.033046136680285, 1.040286556610491, 1.030518454908829, 1.030295809169676, 1.029226987229062, 1.036426538571602, 1.043641727248607, 1.014141637444132),
    @(25, 1.02, 1.02555525487941, 1.027126269869887, 1.034460760183412, 1.041980660673722, 1.030718106631404, 1.031210450974132, 1.030206310070455, 1.037517504273741, 1.045013891748725, 1.014443777279167)
)

foreach ($entry in $data) {
    $r = $entry[0]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $c = $cols[$i]
        $v = $entry[$i + 1]
        $ws.Cells.Item($r, $c).Value = $v
    }
}

Write-Host "Updated $($data.Count) rows (380 kV case)"